# Updated symbol list (coin names/links reshuffled; prices & 1h volume % refreshed).
# Cell-by-cell values taken from the target OOXML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cells: coin name (B) and coinranking link (C) ---
$textUpdates = @(
    @("B9", "WazirX"),
    @("C9", "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"),
    @("B10", "LiechtensteinCryptoassetsExchange"),
    @("C10", "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"),
    @("B11", "MandalaExchangeToken"),
    @("C11", "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"),
    @("B12", "BitrueCoin"),
    @("C12", "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"),
    @("B13", "BitMartToken"),
    @("C13", "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"),
    @("B14", "BitForexToken"),
    @("C14", "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"),
    @("B15", "One"),
    @("C15", "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"),
    @("B16", "TigerCash"),
    @("C16", "https://coinranking.com/coin/6hIn06L2+tigercash-tch"),
    @("B17", "LEO"),
    @("C17", "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"),
    @("B18", "GateToken"),
    @("C18", "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"),
    @("B19", "BTSEToken"),
    @("C19", "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"),
    @("B20", "BitpandaEcosystemToken"),
    @("C20", "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best")
)

foreach ($pair in $textUpdates) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# --- Numeric-looking text cells: price (D) and 1h volume % (E) ---
# These are stored as text in the source sheet (e.g. "263.03", "0.84%"),
# so force Text number format first to stop Excel from coercing the
# assigned string into a float / percentage value, then drop back to the
# workbook-default "Normal" style so no extra formatting is introduced.
$numericTextUpdates = @(
    @("D2", "263.03"),
    @("E2", "0.84%"),
    @("D3", "26.64"),
    @("E3", "-1.90%"),
    @("D4", "4.692"),
    @("E4", "0.42%"),
    @("D5", "0.06100"),
    @("E5", "-1.25%"),
    @("D6", "6.707"),
    @("E6", "0.51%"),
    @("D7", "0.8514"),
    @("E7", "-0.05%"),
    @("D8", "0.9112"),
    @("E8", "-0.60%"),
    @("D9", "0.1411"),
    @("E9", "-0.03%"),
    @("D10", "0.04768"),
    @("E10", "-1.06%"),
    @("D11", "0.07094"),
    @("E11", "0.11%"),
    @("D12", "0.03134"),
    @("E12", "0.77%"),
    @("D13", "0.09048"),
    @("E13", "-0.01%"),
    @("D14", "0.001528"),
    @("E14", "-0.71%"),
    @("D15", "0.0006191"),
    @("E15", "0.09%"),
    @("D16", "0.005990"),
    @("E16", "-1.69%"),
    @("D17", "3.452"),
    @("E17", "0.05%"),
    @("D18", "3.165"),
    @("E18", "0.42%"),
    @("D19", "2.146"),
    @("E19", "-1.54%"),
    @("D20", "0.3072"),
    @("E20", "-0.20%"),
    @("D22", "4.131"),
    @("E22", "1.17%"),
    @("D23", "0.04229"),
    @("E23", "-0.16%"),
    @("D24", "0.001178"),
    @("E24", "-3.07%"),
    @("D25", "0.004066"),
    @("E25", "6.93%"),
    @("E26", "0.02%"),
    @("D40", "0.03931"),
    @("E40", "1.41%"),
    @("D41", "0.1115"),
    @("E41", "0.12%"),
    @("D42", "0.004174"),
    @("E42", "1.98%"),
    @("D43", "0.002110"),
    @("E43", "-3.92%"),
    @("D44", "0.01149"),
    @("E44", "-29.63%"),
    @("D45", "0.00005082"),
    @("E45", "-1.38%"),
    @("E46", "0.01%"),
    @("E48", "59.09%"),
    @("E49", "0.01%"),
    @("E50", "0.01%")
)

foreach ($pair in $numericTextUpdates) {
    $cell = $ws.Range($pair[0])
    $cell.NumberFormat = "@"
    $cell.Value = $pair[1]
    $cell.Style = "Normal"
}

